$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "Emmanuel"
    3  = "Nadiya"
    4  = "Moussa"
    5  = "Hajer"
    6  = "Charly"
    7  = "Floriane"
    8  = "Dragos"
    9  = "Sofia"
    10 = "Mengstu"
    11 = "Choti"
    12 = "Preeti"
    13 = "Klebert"
    14 = "Kenny"
    15 = "Hanieh"
    16 = "Augustin"
    17 = "Yassine"
    18 = "Caterina"
    19 = "Evi"
    20 = "Megan"
    21 = "Fang"
    22 = "Aida"
    23 = "Marc2"
    24 = "Younes"
    25 = "Jordi"
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
